# Apply the change: add a new "to-do" item to the "coisas a fazer" sheet
# describing implementing Flask-WTF after the 1st version.

$wb = $excel.ActiveWorkbook

# The target worksheet is "coisas a fazer" (the 2nd sheet)
$ws = $wb.Worksheets.Item("coisas a fazer")

# Row 14: number (A14) copies the style/formatting pattern of the numbering
# column (same as the other numbered rows), and the description text (B14)
# uses the same style as the other plain rows (e.g. B9:B12).
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B12").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "Implementar o Flask-WTF, mas só depois da 1ºVersão"

# Update the selection to match the recorded state after the edit
$ws.Range("B15").Select() | Out-Null
